$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$ws.Range("A2").Value = "http://172.16.2.61:1616/UI#"
$ws.Range("B2").Value = "Administrator"
$ws.Range("C2").Value = "Tetherfi@930"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Tetherfi@930")
$ws.Range("C2").VerticalAlignment = -4108

$ws.Range("E10").Select()
